$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'91.244.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.87%  "
$ws.Range("D3").Value = "'3.183.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'216.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.12%  "
$ws.Range("D6").Value = "'634.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.94%  "
$ws.Range("D7").Value = "'0.394"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.31%  "
$ws.Range("D8").Value = "'0.717"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.06%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "'3.180.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.27%  "
$ws.Range("D11").Value = "'0.568"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.67%  "
$ws.Range("E12").Value = "  +3.17%  "
$ws.Range("E13").Value = "  +3.79%  "
$ws.Range("D14").Value = "'90.806.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.52%  "
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("D16").Value = "'3.769.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "'32.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").Value = "'3.194.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("D19").Value = "'3.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.19%  "
$ws.Range("D20").Value = "'0.0000215"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +57.14%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "'13.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.09%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "'435.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.59%  "
$ws.Range("D23").Value = "'8.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").Value = "'4.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.83%  "
$ws.Range("D25").Value = "'5.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("E26").Value = "  -5.14%  "
$ws.Range("D27").Value = "'80.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.14%  "
$ws.Range("D28").Value = "'3.352.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.35%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  -2.04%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  +29.85%  "
$ws.Range("D33").Value = "'8.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.07%  "
$ws.Range("D34").Value = "'516.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.84%  "
$ws.Range("D35").Value = "'6.98"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.66%  "
$ws.Range("E36").Value = "  +1.74%  "
$ws.Range("D37").Value = "'1.30"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.38%  "
$ws.Range("D38").Value = "'22.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.07%  "
$ws.Range("D39").Value = "'22.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.81%  "
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").Value = "'0.126"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.72%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("E44").Value = "  -1.11%  "
$ws.Range("D45").Value = "'147.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.81%  "
$ws.Range("D46").Value = "'44.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.10%  "
$ws.Range("D47").Value = "'169.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.36%  "
$ws.Range("D48").Value = "'0.125"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").Value = "'0.739"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.00%  "
$ws.Range("D50").Value = "'24.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.09%  "
$ws.Range("D51").Value = "'0.608"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.96%  "
